$wb = $excel.ActiveWorkbook

# --- Update the publication Date on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-02-21T11:59:56+00:00"

# --- Add a new concept row (OTHER / Undefined Data Type) on the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Write the new values first. Force column A ("Level") to be stored as text
# (matching the existing "1" values in that column, which are shared strings,
# not numbers).
$concepts.Range("A9").NumberFormat = "@"
$concepts.Range("A9").Value = "1"
$concepts.Range("B9").Value = "OTHER"
$concepts.Range("C9").Value = "Undefined Data Type"

# Now copy the formatting from the row above so the new row matches the
# existing table styling exactly (border/fill/font/alignment).
$concepts.Range("A8:D8").Copy()
$concepts.Range("A9:D9").PasteSpecial(-4122)
